$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '36.528.45'
$ws.Range('E2').Value = '  +0.21%  '
$ws.Range('D3').Value = '1.952.83'
$ws.Range('E3').Value = '  +0.89%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = '243.27'
$ws.Range('E5').Value = '  +0.52%  '
$ws.Range('D6').Value = '0.623'
$ws.Range('E6').Value = '  +2.38%  '
$ws.Range('D7').Value = '59.78'
$ws.Range('E7').Value = '  +6.16%  '
$ws.Range('E8').Value = '  -0.02%  '
$ws.Range('D9').Value = '0.378'
$ws.Range('E9').Value = '  +5.46%  '
$ws.Range('D10').Value = '0.0788'
$ws.Range('E10').Value = '  -2.92%  '
$ws.Range('E11').Value = '  +0.74%  '
$ws.Range('D12').Value = '14.12'
$ws.Range('E12').Value = '  +6.77%  '
$ws.Range('D13').Value = '2.241.65'
$ws.Range('E13').Value = '  +0.86%  '
$ws.Range('E14').Value = '  +3.01%  '
$ws.Range('D15').Value = '21.48'
$ws.Range('E15').Value = '  +2.13%  '
$ws.Range('E16').Value = '  +2.46%  '
$ws.Range('D17').Value = '1.950.26'
$ws.Range('E17').Value = '  +0.58%  '
$ws.Range('D18').Value = '36.459.18'
$ws.Range('E18').Value = '  +0.27%  '
$ws.Range('D19').Value = '69.15'
$ws.Range('E19').Value = '  +0.35%  '
$ws.Range('D20').Value = '0.0₃0851'
$ws.Range('E20').Value = '  -0.06%  '
$ws.Range('D21').Value = '229.15'
$ws.Range('E21').Value = '  +1.19%  '
$ws.Range('D22').Value = '5.05'
$ws.Range('E22').Value = '  +2.43%  '
$ws.Range('E23').Value = '  +0.06%  '
$ws.Range('E24').Value = '  +2.58%  '
$ws.Range('D25').Value = '2.37'
$ws.Range('E25').Value = '  +3.47%  '
$ws.Range('E26').Value = '  +8.15%  '
$ws.Range('D27').Value = '9.12'
$ws.Range('E27').Value = '  +0.75%  '
$ws.Range('D28').Value = '160.35'
$ws.Range('E28').Value = '  +0.74%  '
$ws.Range('D29').Value = '19.22'
$ws.Range('E29').Value = '  +0.97%  '
$ws.Range('D30').Value = '1.30'
$ws.Range('E30').Value = '  +20.40%  '
$ws.Range('E31').Value = '  +2.17%  '
$ws.Range('D32').Value = '4.75'
$ws.Range('E32').Value = '  +4.60%  '
$ws.Range('D33').Value = '0.0609'
$ws.Range('E33').Value = '  +0.00%  '
$ws.Range('D34').Value = '4.43'
$ws.Range('E34').Value = '  +7.85%  '
$ws.Range('B35').Value = 'BinanceUSD'
$ws.Range('C35').Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  -0.07%  '
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D36').Value = '3.43'
$ws.Range('E36').Value = '  +10.30%  '
$ws.Range('B37').Value = 'LidoDAOToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D37').Value = '2.26'
$ws.Range('E37').Value = '  +4.63%  '
$ws.Range('E38').Value = '  -0.50%  '
$ws.Range('D39').Value = '5.43'
$ws.Range('E39').Value = '  -10.70%  '
$ws.Range('D40').Value = '0.0965'
$ws.Range('E40').Value = '  -1.85%  '
$ws.Range('E41').Value = '  +0.17%  '
$ws.Range('E42').Value = '  +2.15%  '
$ws.Range('E43').Value = '  +1.15%  '
$ws.Range('D44').Value = '15.74'
$ws.Range('E44').Value = '  +0.93%  '
$ws.Range('D45').Value = '1.360.17'
$ws.Range('E45').Value = '  +2.34%  '
$ws.Range('D46').Value = '88.66'
$ws.Range('E46').Value = '  +3.76%  '
$ws.Range('D47').Value = '1.02'
$ws.Range('E47').Value = '  +0.57%  '
$ws.Range('D48').Value = '7.14'
$ws.Range('E48').Value = '  +0.70%  '
$ws.Range('E49').Value = '  +0.94%  '
$ws.Range('D50').Value = '45.29'
$ws.Range('E50').Value = '  +5.77%  '
$ws.Range('D51').Value = '2.137.32'
$ws.Range('E51').Value = '  +1.12%  '
